$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $savedStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $savedStyle
}

# Row 2
$ws.Range("D2").Value = "68.049.86"
$ws.Range("E2").Value = "  +1.38%  "

# Row 3
$ws.Range("D3").Value = "2.515.12"
$ws.Range("E3").Value = "  +1.14%  "

# Row 4
Set-TextValue $ws.Range("D4") "0.998"
$ws.Range("E4").Value = "  -0.18%  "

# Row 5
Set-TextValue $ws.Range("D5") "589.68"
$ws.Range("E5").Value = "  +1.02%  "

# Row 6
Set-TextValue $ws.Range("D6") "177.91"
$ws.Range("E6").Value = "  +3.61%  "

# Row 7
$ws.Range("E7").Value = "  -0.08%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.516"
$ws.Range("E8").Value = "  +0.59%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.142"
$ws.Range("E9").Value = "  +3.58%  "

# Row 10
$ws.Range("E10").Value = "  -0.28%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.341"
$ws.Range("E11").Value = "  +2.62%  "

# Row 12
$ws.Range("E12").Value = "  +0.64%  "

# Row 13
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D13") "25.78"
$ws.Range("E13").Value = "  +1.61%  "

# Row 14
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.926.11"
$ws.Range("E14").Value = "  -0.51%  "

# Row 15
$ws.Range("D15").Value = "67.879.85"
$ws.Range("E15").Value = "  +1.19%  "

# Row 16
$ws.Range("E16").Value = "  +1.18%  "

# Row 17
$ws.Range("D17").Value = "2.493.56"
$ws.Range("E17").Value = "  -1.60%  "

# Row 18
Set-TextValue $ws.Range("D18") "11.04"
$ws.Range("E18").Value = "  +0.35%  "

# Row 19
Set-TextValue $ws.Range("D19") "7.55"
$ws.Range("E19").Value = "  +1.96%  "

# Row 20
Set-TextValue $ws.Range("D20") "353.67"
$ws.Range("E20").Value = "  +1.40%  "

# Row 21
Set-TextValue $ws.Range("D21") "4.12"

# Row 22
$ws.Range("E22").Value = "  +0.10%  "

# Row 23
Set-TextValue $ws.Range("D23") "70.91"
$ws.Range("E23").Value = "  +3.42%  "

# Row 24
Set-TextValue $ws.Range("D24") "4.35"
$ws.Range("E24").Value = "  +3.15%  "

# Row 25
Set-TextValue $ws.Range("D25") "1.77"
$ws.Range("E25").Value = "  -0.74%  "

# Row 26
Set-TextValue $ws.Range("D26") "9.16"
$ws.Range("E26").Value = "  -1.10%  "

# Row 27
$ws.Range("D27").Value = "2.592.18"
$ws.Range("E27").Value = "  -0.88%  "

# Row 28
Set-TextValue $ws.Range("D28") "0.992"
$ws.Range("E28").Value = "  -0.65%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0923"
$ws.Range("E29").Value = "  +2.32%  "

# Row 30
Set-TextValue $ws.Range("D30") "509.22"
$ws.Range("E30").Value = "  +0.00%  "

# Row 31
Set-TextValue $ws.Range("D31") "7.89"
$ws.Range("E31").Value = "  +1.68%  "

# Row 32
$ws.Range("E32").Value = "  +2.96%  "

# Row 33
$ws.Range("E33").Value = "  +0.98%  "

# Row 34
$ws.Range("E34").Value = "  -0.10%  "

# Row 35
$ws.Range("E35").Value = "  +4.51%  "

# Row 36
Set-TextValue $ws.Range("D36") "164.70"
$ws.Range("E36").Value = "  +3.13%  "

# Row 37
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D37") "18.44"
$ws.Range("E37").Value = "  +1.23%  "

# Row 38
$ws.Range("B38").Value = "WhiteBITCoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue $ws.Range("D38") "18.66"
$ws.Range("E38").Value = "  -0.21%  "

# Row 39
$ws.Range("E39").Value = "  +0.57%  "

# Row 40
$ws.Range("E40").Value = "  +0.06%  "

# Row 41
$ws.Range("E41").Value = "  +3.26%  "

# Row 42
Set-TextValue $ws.Range("D42") "4.91"
$ws.Range("E42").Value = "  +2.07%  "

# Row 43
Set-TextValue $ws.Range("D43") "0.331"
$ws.Range("E43").Value = "  +0.80%  "

# Row 44
Set-TextValue $ws.Range("D44") "2.49"
$ws.Range("E44").Value = "  +5.70%  "

# Row 45
Set-TextValue $ws.Range("D45") "145.76"
$ws.Range("E45").Value = "  +2.42%  "

# Row 46
Set-TextValue $ws.Range("D46") "3.55"
$ws.Range("E46").Value = "  +3.04%  "

# Row 47
Set-TextValue $ws.Range("D47") "0.521"
$ws.Range("E47").Value = "  +1.54%  "

# Row 48
$ws.Range("D48").Value = "0.0₆0259"
$ws.Range("E48").Value = "  +3.87%  "

# Row 49
$ws.Range("E49").Value = "  +1.58%  "

# Row 50
Set-TextValue $ws.Range("D50") "1.60"
$ws.Range("E50").Value = "  +1.98%  "

# Row 51
Set-TextValue $ws.Range("D51") "0.589"
$ws.Range("E51").Value = "  +1.04%  "
